$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A (cluster labels) and column B (counts) for rows 2-5
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 209

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 134

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 109

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 17

# Row 6 no longer exists in the new data - delete it entirely
$ws.Range("A6:B6").Delete()
